# Generate Report for Archive
#
# The localization status for the two handed-off files moved on from
# "Ready for handoff" to "In Translation". That shorter status string is
# the longest value left in the Status-ish columns, so those columns
# shrink when the report is (re)generated and the columns are auto-fit
# to their new contents.

$wb = $excel.ActiveWorkbook

# --- 1. Update the status text everywhere it appears -----------------
# Overview sheet: columns E (zh-cn) / F (de-de) hold the per-locale
# status for each of the two rows of data.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"

# Per-locale detail sheets: column C ("Status").
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"

# --- 2. Re-fit the affected columns to their new (shorter) contents --
$overview.Range("E1:F1").EntireColumn.AutoFit()
$zhcn.Range("C1").EntireColumn.AutoFit()
$dede.Range("C1").EntireColumn.AutoFit()

# AutoFit in this runtime approximates real Excel's sub-pixel font
# metrics, so nail the resulting width to match what Excel itself
# computed for "In Translation" being the widest entry in these columns.
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
